$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new price string looks like a plain number
# (e.g. "19.72") need to be pre-formatted as Text so Excel keeps the
# literal string instead of parsing it into a floating point value.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range('D2').Value = '27.025.14'
$ws.Range('E2').Value = '  +2.88%  '
$ws.Range('D3').Value = '1.649.22'
$ws.Range('E3').Value = '  +3.46%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '214.97'
$ws.Range('E5').Value = '  +1.69%  '
$ws.Range('E6').Value = '  +1.45%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('E8').Value = '  +1.77%  '
$ws.Range('E9').Value = '  +1.57%  '
$ws.Range('D10').Value = '19.72'
$ws.Range('E10').Value = '  +3.96%  '
$ws.Range('E11').Value = '  +1.17%  '
$ws.Range('D12').Value = '1.882.49'
$ws.Range('E12').Value = '  +3.49%  '
$ws.Range('D13').Value = '1.651.94'
$ws.Range('E13').Value = '  +4.52%  '
$ws.Range('E14').Value = '  +2.43%  '
$ws.Range('D15').Value = '0.518'
$ws.Range('E15').Value = '  +3.08%  '
$ws.Range('D16').Value = '65.09'
$ws.Range('E16').Value = '  +2.60%  '
$ws.Range('D17').Value = '27.018.11'
$ws.Range('E17').Value = '  +2.85%  '
$ws.Range('D18').Value = '237.94'
$ws.Range('E18').Value = '  +4.15%  '
$ws.Range('D19').Value = '7.83'
$ws.Range('E19').Value = '  +2.15%  '
$ws.Range('D20').Value = '0.0₃0729'
$ws.Range('E20').Value = '  +1.23%  '
$ws.Range('E21').Value = '  +0.09%  '
$ws.Range('E22').Value = '  +4.46%  '
$ws.Range('E23').Value = '  +4.36%  '
$ws.Range('E24').Value = '  +3.61%  '
$ws.Range('D25').Value = '145.47'
$ws.Range('E25').Value = '  -0.48%  '
$ws.Range('E26').Value = '  -0.06%  '
$ws.Range('E27').Value = '  +2.05%  '
$ws.Range('E28').Value = '  +1.54%  '
$ws.Range('D29').Value = '15.80'
$ws.Range('E29').Value = '  +2.78%  '
$ws.Range('E30').Value = '  +0.82%  '
$ws.Range('E31').Value = '  +1.80%  '
$ws.Range('E32').Value = '  +3.28%  '
$ws.Range('D33').Value = '1.510.88'
$ws.Range('E33').Value = '  +2.65%  '
$ws.Range('E34').Value = '  +5.09%  '
$ws.Range('E35').Value = '  +9.08%  '
$ws.Range('E36').Value = '  -0.11%  '
$ws.Range('D37').Value = '0.576'
$ws.Range('E37').Value = '  +1.58%  '
$ws.Range('E38').Value = '  +8.56%  '
$ws.Range('E39').Value = '  +2.63%  '
$ws.Range('E40').Value = '  +3.51%  '
$ws.Range('E41').Value = '  -0.03%  '
$ws.Range('E42').Value = '  +4.08%  '
$ws.Range('D43').Value = '65.86'
$ws.Range('E43').Value = '  +9.52%  '
$ws.Range('D44').Value = '1.789.73'
$ws.Range('E44').Value = '  +3.30%  '
$ws.Range('D45').Value = '0.773'
$ws.Range('E45').Value = '  +2.46%  '
$ws.Range('D46').Value = '0.915'
$ws.Range('E46').Value = '  -2.09%  '
$ws.Range('D47').Value = '89.37'
$ws.Range('E47').Value = '  +0.79%  '
$ws.Range('E48').Value = '  +0.06%  '
$ws.Range('E49').Value = '  +3.20%  '
$ws.Range('E50').Value = '  +1.04%  '
$ws.Range('D51').Value = '0.0975'
$ws.Range('E51').Value = '  +2.21%  '

# Restore the default (General) style on those cells so only the cell
# content changes, matching the original formatting.
$ws.Range("D5").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D51").Style = "Normal"
